$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: advance the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the three price cells from 94.3 to 185
$ws.Range("D27").Value = 185
$ws.Range("D28").Value = 185
$ws.Range("D29").Value = 185
